# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de):
#  - Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#  - Two new report columns are populated:
#       F = Latest Target File    (same source .md file as column A)
#       G = Latest Handback File  (the handed-back .xlf file, mirrors column D)
#  - Latest Handback DateTime (H) is stamped with the real handback time
#
# The Overview sheet's per-language status cells share the same text, so they
# flip to the new wording too.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR packing of RGB(0x64,0x95,0xED) -> matches the workbook's HyperLink font color FF6495ED

# ---------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells (same wording
# as the per-language sheets' Status column, so it must match).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zhTargetName = "e6f30ad4-923a-4500-9d07-14030507f370.md"
$zhTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5f26bbc8c24dabe60a6e0793c43c22f6928cd7bd/e2e/e6f30ad4-923a-4500-9d07-14030507f370.md"
$zhSourceAName = "ffff95e8289a-0bb2-4d81-a22f-748725d7e54d.md"
$zhSourceAUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5f26bbc8c24dabe60a6e0793c43c22f6928cd7bd/e2e/ffff95e8289a-0bb2-4d81-a22f-748725d7e54d.md"
$zhHandbackName = "e6f30ad4-923a-4500-9d07-14030507f370.bd19de4e2d5fad4ee8228e061147df4a76f02433.zh-cn.xlf"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b39c37569c89cfea8801ad0e8f44db33d3bcccb6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e6f30ad4-923a-4500-9d07-14030507f370.bd19de4e2d5fad4ee8228e061147df4a76f02433.zh-cn.xlf"

# Re-create the existing A2/D2/A3/D3 hyperlinks alongside the two new
# columns so every hyperlink on the sheet is (re)inserted in natural
# left-to-right, top-to-bottom report order (rId2..rId9).
$zh.Range("A2").Hyperlinks.Delete()
$zh.Range("D2").Hyperlinks.Delete()
$zh.Range("A3").Hyperlinks.Delete()
$zh.Range("D3").Hyperlinks.Delete()

$zh.Hyperlinks.Add($zh.Range("A2"), $zhTargetUrl, "", "", $zhTargetName)
$zh.Hyperlinks.Add($zh.Range("D2"), $zhHandbackUrl, "", "", $zhHandbackName)
$zh.Hyperlinks.Add($zh.Range("F2"), $zhTargetUrl, "", "", $zhTargetName)
$zh.Hyperlinks.Add($zh.Range("G2"), $zhHandbackUrl, "", "", $zhHandbackName)
$zh.Hyperlinks.Add($zh.Range("A3"), $zhSourceAUrl, "", "", $zhSourceAName)
$zh.Hyperlinks.Add($zh.Range("D3"), $zhHandbackUrl, "", "", $zhHandbackName)
$zh.Hyperlinks.Add($zh.Range("F3"), $zhTargetUrl, "", "", $zhTargetName)
$zh.Hyperlinks.Add($zh.Range("G3"), $zhHandbackUrl, "", "", $zhHandbackName)

$zh.Range("F2:G3").Font.Underline = 2
$zh.Range("F2:G3").Font.Color = $hyperlinkColor

$zh.Range("H2").Value = "2016-03-24 23:16:07"
$zh.Range("H3").Value = "2016-03-24 23:16:07"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$deTargetName = "e6f30ad4-923a-4500-9d07-14030507f370.md"
$deTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5f26bbc8c24dabe60a6e0793c43c22f6928cd7bd/e2e/e6f30ad4-923a-4500-9d07-14030507f370.md"
$deSourceAName = "ffff95e8289a-0bb2-4d81-a22f-748725d7e54d.md"
$deSourceAUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5f26bbc8c24dabe60a6e0793c43c22f6928cd7bd/e2e/ffff95e8289a-0bb2-4d81-a22f-748725d7e54d.md"
$deHandbackName = "e6f30ad4-923a-4500-9d07-14030507f370.bd19de4e2d5fad4ee8228e061147df4a76f02433.de-de.xlf"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6933a812214763bb877b20d34c0f946b0687a99c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e6f30ad4-923a-4500-9d07-14030507f370.bd19de4e2d5fad4ee8228e061147df4a76f02433.de-de.xlf"

$de.Range("A2").Hyperlinks.Delete()
$de.Range("D2").Hyperlinks.Delete()
$de.Range("A3").Hyperlinks.Delete()
$de.Range("D3").Hyperlinks.Delete()

$de.Hyperlinks.Add($de.Range("A2"), $deTargetUrl, "", "", $deTargetName)
$de.Hyperlinks.Add($de.Range("D2"), $deHandbackUrl, "", "", $deHandbackName)
$de.Hyperlinks.Add($de.Range("F2"), $deTargetUrl, "", "", $deTargetName)
$de.Hyperlinks.Add($de.Range("G2"), $deHandbackUrl, "", "", $deHandbackName)
$de.Hyperlinks.Add($de.Range("A3"), $deSourceAUrl, "", "", $deSourceAName)
$de.Hyperlinks.Add($de.Range("D3"), $deHandbackUrl, "", "", $deHandbackName)
$de.Hyperlinks.Add($de.Range("F3"), $deTargetUrl, "", "", $deTargetName)
$de.Hyperlinks.Add($de.Range("G3"), $deHandbackUrl, "", "", $deHandbackName)

$de.Range("F2:G3").Font.Underline = 2
$de.Range("F2:G3").Font.Color = $hyperlinkColor

$de.Range("H2").Value = "2016-03-24 23:16:14"
$de.Range("H3").Value = "2016-03-24 23:16:14"
